$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet (3rd sheet)
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# trailing columns (Late / heading spacer / Outstanding) one column to
# the right. The newly inserted column inherits the width that the
# "In Advance" column (M) already has, matching what Excel does when a
# column is inserted next to existing, explicitly-sized columns.
$mColumnWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mColumnWidth

# Make "Repayment schedule" the active sheet and move the selection to
# the single cell R6, as captured by the workbook view state.
$ws.Activate()
$ws.Range("R6").Select()
